# Update building block types for SRA - Sequencing template

$wb = $excel.ActiveWorkbook

$wsMeta  = $wb.Worksheets.Item("isa_template")
$wsTable = $wb.Worksheets.Item("Next generation sequencing")

# Bump template version 1.0.1 -> 1.0.2
$wsMeta.Range("B4").Value = "1.0.2"

# Rename building block headers (row 1) - table headers stay in sync automatically
$wsTable.Range("R1").Value  = "Characteristic [library source]"
$wsTable.Range("AD1").Value = "Component [next generation sequencing instrument model]"
$wsTable.Range("AJ1").Value = "Output [Data]"

# Update ontology term URLs / source refs in the example data row (row 2)
$wsTable.Range("Q2").Value  = "https://bioregistry.io/EFO:0008896"
$wsTable.Range("W2").Value  = "https://bioregistry.io/NCIT:C28408"
$wsTable.Range("Z2").Value  = "http://purl.org/nfdi4plants/ontology/dpbo/DPBO_0000086"
$wsTable.Range("AE2").Value = "EFO"
$wsTable.Range("AF2").Value = "https://bioregistry.io/EFO:0008563"
